$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI re-run ("Natmi following Dr Hou advice"): the ligand- and receptor-
# expressing cell counts (columns E and K) increase from 1 to 3 for every data
# row, which cascades into the average/total expression, derived-specificity,
# and edge-weight columns (G,H,I,J,M,N,O,P,Q,R,S,T). Columns A-D, F and L are
# unaffected, so only the recomputed cells are written back explicitly.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.191447666666667
$ws.Range("H2").Value = 24.574343
$ws.Range("I2").Value = 0.185794284429433
$ws.Range("J2").Value = 0.185794284429433
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.021894333333333
$ws.Range("N2").Value = 9.065683
$ws.Range("O2").Value = 0.1464771679819186
$ws.Range("P2").Value = 0.1464771679819185
$ws.Range("Q2").Value = 24.75368928569656
$ws.Range("R2").Value = 222.783203571269
$ws.Range("S2").Value = 0.02721462061045041
$ws.Range("T2").Value = 0.0272146206104504

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.191447666666667
$ws.Range("H3").Value = 24.574343
$ws.Range("I3").Value = 0.185794284429433
$ws.Range("J3").Value = 0.185794284429433
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.799695333333333
$ws.Range("N3").Value = 17.399086
$ws.Range("O3").Value = 0.2811226515149324
$ws.Range("P3").Value = 0.2811226515149324
$ws.Range("Q3").Value = 47.50790080561089
$ws.Range("R3").Value = 427.571107250498
$ws.Range("S3").Value = 0.05223098187512172
$ws.Range("T3").Value = 0.05223098187512172

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.191447666666667
$ws.Range("H4").Value = 24.574343
$ws.Range("I4").Value = 0.185794284429433
$ws.Range("J4").Value = 0.185794284429433
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.259039333333333
$ws.Range("N4").Value = 18.777118
$ws.Range("O4").Value = 0.303387959572633
$ws.Range("P4").Value = 0.303387959572633
$ws.Range("Q4").Value = 51.27059314260822
$ws.Range("R4").Value = 461.4353382834739
$ws.Range("S4").Value = 0.05636774885330308
$ws.Range("T4").Value = 0.05636774885330308

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.191447666666667
$ws.Range("H5").Value = 24.574343
$ws.Range("I5").Value = 0.185794284429433
$ws.Range("J5").Value = 0.185794284429433
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.549851333333334
$ws.Range("N5").Value = 16.649554
$ws.Range("O5").Value = 0.2690122209305161
$ws.Range("P5").Value = 0.2690122209305161
$ws.Range("Q5").Value = 45.46131675478023
$ws.Range("R5").Value = 409.151850793022
$ws.Range("S5").Value = 0.04998093309055777
$ws.Range("T5").Value = 0.04998093309055776

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.317702
$ws.Range("H6").Value = 51.95310600000001
$ws.Range("I6").Value = 0.3927913821808575
$ws.Range("J6").Value = 0.3927913821808576
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.021894333333333
$ws.Range("N6").Value = 9.065683
$ws.Range("O6").Value = 0.1464771679819186
$ws.Range("P6").Value = 0.1464771679819185
$ws.Range("Q6").Value = 52.33226554015534
$ws.Range("R6").Value = 470.9903898613981
$ws.Range("S6").Value = 0.05753496926955544
$ws.Range("T6").Value = 0.05753496926955544

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.317702
$ws.Range("H7").Value = 51.95310600000001
$ws.Range("I7").Value = 0.3927913821808575
$ws.Range("J7").Value = 0.3927913821808576
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.799695333333333
$ws.Range("N7").Value = 17.399086
$ws.Range("O7").Value = 0.2811226515149324
$ws.Range("P7").Value = 0.2811226515149324
$ws.Range("Q7").Value = 100.4373954734573
$ws.Range("R7").Value = 903.9365592611161
$ws.Range("S7").Value = 0.1104225548508978
$ws.Range("T7").Value = 0.1104225548508979

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.317702
$ws.Range("H8").Value = 51.95310600000001
$ws.Range("I8").Value = 0.3927913821808575
$ws.Range("J8").Value = 0.3927913821808576
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.259039333333333
$ws.Range("N8").Value = 18.777118
$ws.Range("O8").Value = 0.303387959572633
$ws.Range("P8").Value = 0.303387959572633
$ws.Range("Q8").Value = 108.3921779809453
$ws.Range("R8").Value = 975.529601828508
$ws.Range("S8").Value = 0.1191681759775646
$ws.Range("T8").Value = 0.1191681759775646

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.317702
$ws.Range("H9").Value = 51.95310600000001
$ws.Range("I9").Value = 0.3927913821808575
$ws.Range("J9").Value = 0.3927913821808576
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.549851333333334
$ws.Range("N9").Value = 16.649554
$ws.Range("O9").Value = 0.2690122209305161
$ws.Range("P9").Value = 0.2690122209305161
$ws.Range("Q9").Value = 96.11067153496934
$ws.Range("R9").Value = 864.9960438147242
$ws.Range("S9").Value = 0.1056656820828396
$ws.Range("T9").Value = 0.1056656820828396

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.35128266666667
$ws.Range("H10").Value = 37.053848
$ws.Range("I10").Value = 0.2801455637905346
$ws.Range("J10").Value = 0.2801455637905346
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.021894333333333
$ws.Range("N10").Value = 9.065683
$ws.Range("O10").Value = 0.1464771679819186
$ws.Range("P10").Value = 0.1464771679819185
$ws.Range("Q10").Value = 37.32427109979822
$ws.Range("R10").Value = 335.918439898184
$ws.Range("S10").Value = 0.04103492880673541
$ws.Range("T10").Value = 0.04103492880673541

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 12.35128266666667
$ws.Range("H11").Value = 37.053848
$ws.Range("I11").Value = 0.2801455637905346
$ws.Range("J11").Value = 0.2801455637905346
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.799695333333333
$ws.Range("N11").Value = 17.399086
$ws.Range("O11").Value = 0.2811226515149324
$ws.Range("P11").Value = 0.2811226515149324
$ws.Range("Q11").Value = 71.63367644254757
$ws.Range("R11").Value = 644.7030879829281
$ws.Range("S11").Value = 0.07875526370294073
$ws.Range("T11").Value = 0.07875526370294073

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 12.35128266666667
$ws.Range("H12").Value = 37.053848
$ws.Range("I12").Value = 0.2801455637905346
$ws.Range("J12").Value = 0.2801455637905346
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.259039333333333
$ws.Range("N12").Value = 18.777118
$ws.Range("O12").Value = 0.303387959572633
$ws.Range("P12").Value = 0.303387959572633
$ws.Range("Q12").Value = 77.30716402778489
$ws.Range("R12").Value = 695.7644762500639
$ws.Range("S12").Value = 0.08499279098173518
$ws.Range("T12").Value = 0.08499279098173518

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 12.35128266666667
$ws.Range("H13").Value = 37.053848
$ws.Range("I13").Value = 0.2801455637905346
$ws.Range("J13").Value = 0.2801455637905346
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.549851333333334
$ws.Range("N13").Value = 16.649554
$ws.Range("O13").Value = 0.2690122209305161
$ws.Range("P13").Value = 0.2690122209305161
$ws.Range("Q13").Value = 68.5477825759769
$ws.Range("R13").Value = 616.9300431837921
$ws.Range("S13").Value = 0.07536258029912328
$ws.Range("T13").Value = 0.07536258029912327

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.228371
$ws.Range("H14").Value = 18.685113
$ws.Range("I14").Value = 0.1412687695991749
$ws.Range("J14").Value = 0.1412687695991749
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.021894333333333
$ws.Range("N14").Value = 9.065683
$ws.Range("O14").Value = 0.1464771679819186
$ws.Range("P14").Value = 0.1464771679819185
$ws.Range("Q14").Value = 18.82147903079767
$ws.Range("R14").Value = 169.393311277179
$ws.Range("S14").Value = 0.02069264929517729
$ws.Range("T14").Value = 0.02069264929517728

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.228371
$ws.Range("H15").Value = 18.685113
$ws.Range("I15").Value = 0.1412687695991749
$ws.Range("J15").Value = 0.1412687695991749
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.799695333333333
$ws.Range("N15").Value = 17.399086
$ws.Range("O15").Value = 0.2811226515149324
$ws.Range("P15").Value = 0.2811226515149324
$ws.Range("Q15").Value = 36.12265422296866
$ws.Range("R15").Value = 325.103888006718
$ws.Range("S15").Value = 0.03971385108597212
$ws.Range("T15").Value = 0.03971385108597212

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.228371
$ws.Range("H16").Value = 18.685113
$ws.Range("I16").Value = 0.1412687695991749
$ws.Range("J16").Value = 0.1412687695991749
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.259039333333333
$ws.Range("N16").Value = 18.777118
$ws.Range("O16").Value = 0.303387959572633
$ws.Range("P16").Value = 0.303387959572633
$ws.Range("Q16").Value = 38.98361907159266
$ws.Range("R16").Value = 350.852571644334
$ws.Range("S16").Value = 0.04285924376003006
$ws.Range("T16").Value = 0.04285924376003006

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.228371
$ws.Range("H17").Value = 18.685113
$ws.Range("I17").Value = 0.1412687695991749
$ws.Range("J17").Value = 0.1412687695991749
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.549851333333334
$ws.Range("N17").Value = 16.649554
$ws.Range("O17").Value = 0.2690122209305161
$ws.Range("P17").Value = 0.2690122209305161
$ws.Range("Q17").Value = 34.56653309884467
$ws.Range("R17").Value = 311.0987978896021
$ws.Range("S17").Value = 0.03800302545799541
$ws.Range("T17").Value = 0.0380030254579954
